# Sprint 2 backlog review update
#
# Two backlog items are no longer part of the sprint and are removed:
#   - "Manage stock"                               (row 7, Assignee: Luke)
#   - "View record of items tracked by which employee" (row 8, Assignee: Tristen)
# The remaining "Track when things are received and shipped" item (old row 9,
# Assignee: Carson) shifts up to become the new row 7, and all the blank
# styled rows below it shift up by two as well.
#
# Cosmetic changes made while reviewing the sheet: the Assignee/Feature/etc.
# columns were widened and given explicit widths, and the view was re-zoomed
# with a different cell selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two completed/dropped backlog rows (old rows 7 and 8). This
# shifts row 9 ("Track when things are received and shipped") up to row 7,
# and the formatted-but-empty rows 12-17 up to rows 10-15, exactly matching
# the target layout.
$ws.Rows("7:8").Delete()

# Give columns A-D explicit custom widths (set in points; the host stores
# column widths on a coarse grid, so these are the closest attainable
# values to the authored widths of ~22.55 / ~14.89 / ~22.89 / 45 chars).
$ws.Columns("A:A").ColumnWidth = 21.66625
$ws.Columns("B:B").ColumnWidth = 14
$ws.Columns("C:C").ColumnWidth = 22
$ws.Columns("D:D").ColumnWidth = 44.16625

# Update the view: zoom in further and move the selection.
$excel.ActiveWindow.Zoom = 190
$ws.Range("D8").Select() | Out-Null
